$d = $word.ActiveDocument

# Locate the end-of-document block of empty paragraphs that follows the
# Arduino code sample ("Work according to Button Press"). The first
# paragraph of that run must stay untouched; the four new paragraphs of
# commentary are inserted right after it.
$n = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $n; $i++) {
    $cur = $d.Paragraphs.Item($i).Range.Text
    if ($cur -eq "`r" -and $i -lt $n) {
        $next1 = $d.Paragraphs.Item($i + 1).Range.Text
        if ($next1 -eq "`r") {
            $ok = $true
            for ($k = 0; $k -lt 4; $k++) {
                if (($i + $k) -gt $n) { $ok = $false; break }
                if ($d.Paragraphs.Item($i + $k).Range.Text -ne "`r") { $ok = $false; break }
            }
            if ($ok) {
                $targetIndex = $i
                break
            }
        }
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the trailing run of blank paragraphs to anchor the new content."
}

# Paragraph right after the first blank paragraph of the run; this is the
# slot the new content gets inserted in front of.
$anchorPara = $d.Paragraphs.Item($targetIndex + 1)
$r = $anchorPara.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">' +
'<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' +
'<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' +
'</Relationships></pkg:xmlData></pkg:part>' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
'<w:p><w:r><w:t xml:space="preserve">Not working </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>until</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> comes to ready mode.</w:t></w:r></w:p>' +
'<w:p><w:r><w:t xml:space="preserve">Once pressed not working for processing </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>time period</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>. Until ready mode.</w:t></w:r></w:p>' +
'<w:p><w:r><w:t xml:space="preserve">If one button pressed </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>other</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> button </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>not</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> working.</w:t></w:r></w:p>' +
'<w:p><w:r><w:t xml:space="preserve">Indicate those by a light </w:t></w:r></w:p>' +
'<w:p/><w:p/>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $r.InsertXML($xml)
